{"js": "// Update the date line and every \"a\u00f7b=c, d\" division-answer cell in the\n// table to the new values from this revision. Every old value below is\n// unique in the document, so a plain text search + whole-match replace is\n// unambiguous for each pair.\nconst replacements = [\n  [\"2025-09-30 Tuesday\", \"2025-10-01 Wednesday\"],\n  [\"473\u00f74=118, 1\", \"828\u00f73=276, 0\"],\n  [\"234\u00f73=78, 0\", \"269\u00f79=29, 8\"],\n  [\"588\u00f77=84, 0\", \"699\u00f77=99, 6\"],\n  [\"942\u00f75=188, 2\", \"195\u00f74=48, 3\"],\n  [\"134\u00f79=14, 8\", \"263\u00f72=131, 1\"],\n  [\"153\u00f73=51, 0\", \"992\u00f73=330, 2\"],\n  [\"478\u00f74=119, 2\", \"129\u00f76=21, 3\"],\n  [\"909\u00f76=151, 3\", \"325\u00f78=40, 5\"],\n  [\"875\u00f79=97, 2\", \"250\u00f77=35, 5\"],\n  [\"627\u00f75=125, 2\", \"272\u00f74=68, 0\"],\n  [\"599\u00f73=199, 2\", \"436\u00f74=109, 0\"],\n  [\"495\u00f78=61, 7\", \"526\u00f79=58, 4\"],\n  [\"232\u00f77=33, 1\", \"527\u00f76=87, 5\"],\n  [\"804\u00f75=160, 4\", \"598\u00f77=85, 3\"],\n  [\"723\u00f79=80, 3\", \"683\u00f74=170, 3\"],\n  [\"907\u00f76=151, 1\", \"739\u00f77=105, 4\"],\n  [\"909\u00f79=101, 0\", \"654\u00f77=93, 3\"],\n  [\"104\u00f75=20, 4\", \"894\u00f73=298, 0\"],\n  [\"927\u00f76=154, 3\", \"735\u00f79=81, 6\"],\n  [\"193\u00f79=21, 4\", \"107\u00f79=11, 8\"],\n  [\"799\u00f76=133, 1\", \"646\u00f74=161, 2\"],\n  [\"971\u00f76=161, 5\", \"498\u00f75=99, 3\"],\n  [\"668\u00f74=167, 0\", \"221\u00f78=27, 5\"],\n  [\"600\u00f78=75, 0\", \"930\u00f72=465, 0\"],\n  [\"817\u00f75=163, 2\", \"881\u00f73=293, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"a\u00f7b=c, d\" division-answer cell in the\n# table to the new values from this revision. Every old value is unique\n# in the document, so Find/Replace-All per pair is unambiguous and only\n# touches the intended run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-09-30 Tuesday\", \"2025-10-01 Wednesday\"),\n    @(\"473\u00f74=118, 1\", \"828\u00f73=276, 0\"),\n    @(\"234\u00f73=78, 0\", \"269\u00f79=29, 8\"),\n    @(\"588\u00f77=84, 0\", \"699\u00f77=99, 6\"),\n    @(\"942\u00f75=188, 2\", \"195\u00f74=48, 3\"),\n    @(\"134\u00f79=14, 8\", \"263\u00f72=131, 1\"),\n    @(\"153\u00f73=51, 0\", \"992\u00f73=330, 2\"),\n    @(\"478\u00f74=119, 2\", \"129\u00f76=21, 3\"),\n    @(\"909\u00f76=151, 3\", \"325\u00f78=40, 5\"),\n    @(\"875\u00f79=97, 2\", \"250\u00f77=35, 5\"),\n    @(\"627\u00f75=125, 2\", \"272\u00f74=68, 0\"),\n    @(\"599\u00f73=199, 2\", \"436\u00f74=109, 0\"),\n    @(\"495\u00f78=61, 7\", \"526\u00f79=58, 4\"),\n    @(\"232\u00f77=33, 1\", \"527\u00f76=87, 5\"),\n    @(\"804\u00f75=160, 4\", \"598\u00f77=85, 3\"),\n    @(\"723\u00f79=80, 3\", \"683\u00f74=170, 3\"),\n    @(\"907\u00f76=151, 1\", \"739\u00f77=105, 4\"),\n    @(\"909\u00f79=101, 0\", \"654\u00f77=93, 3\"),\n    @(\"104\u00f75=20, 4\", \"894\u00f73=298, 0\"),\n    @(\"927\u00f76=154, 3\", \"735\u00f79=81, 6\"),\n    @(\"193\u00f79=21, 4\", \"107\u00f79=11, 8\"),\n    @(\"799\u00f76=133, 1\", \"646\u00f74=161, 2\"),\n    @(\"971\u00f76=161, 5\", \"498\u00f75=99, 3\"),\n    @(\"668\u00f74=167, 0\", \"221\u00f78=27, 5\"),\n    @(\"600\u00f78=75, 0\", \"930\u00f72=465, 0\"),\n    @(\"817\u00f75=163, 2\", \"881\u00f73=293, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
